$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1923.5
$ws.Range("I40").Value = 1824.8572
$ws.Range("K40").Value = 1824.8572
$ws.Range("M40").Value = -1649.8572
$ws.Range("H74").Value = 5055.8887
$ws.Range("I74").Value = 5325.375
$ws.Range("J74").Value = 2900
$ws.Range("K74").Value = 5325.375
$ws.Range("L74").Value = 2900
$ws.Range("M74").Value = -4389.375
$ws.Range("N74").Value = -4772
$ws.Range("H77").Value = 5055.8887
$ws.Range("I77").Value = 5325.375
$ws.Range("J77").Value = 2900
$ws.Range("K77").Value = 26626.875
$ws.Range("L77").Value = 14500
$ws.Range("M77").Value = -21946.875
$ws.Range("N77").Value = -23860
$ws.Range("H132").Value = 49427964
$ws.Range("I132").Value = 57336148
$ws.Range("J132").Value = 1801.5
$ws.Range("K132").Value = 172008444
$ws.Range("L132").Value = 5404.5
$ws.Range("M132").Value = -172005914
$ws.Range("N132").Value = -10464.5
$ws.Range("H135").Value = 4555.25
$ws.Range("I135").Value = 5732.25
$ws.Range("J135").Value = 2789.75
$ws.Range("K135").Value = 51590.25
$ws.Range("L135").Value = 25107.75
$ws.Range("M135").Value = -49055.25
$ws.Range("N135").Value = -30177.75
$ws.Range("H137").Value = 18740284
$ws.Range("I137").Value = 323043.94
$ws.Range("J137").Value = 58824868
$ws.Range("K137").Value = 969131.8200000001
$ws.Range("L137").Value = 176474604
$ws.Range("M137").Value = -966581.8200000001
$ws.Range("N137").Value = -176479704
$ws.Range("H138").Value = 1714.5294
$ws.Range("I138").Value = 1175.575
$ws.Range("J138").Value = 2193.6
$ws.Range("K138").Value = 3526.725
$ws.Range("L138").Value = 6580.799999999999
$ws.Range("M138").Value = 1613.275
$ws.Range("N138").Value = -16860.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 40400
$ws.Range("J58").Value = 40400
$ws.Range("L58").Value = 40400
$ws.Range("N58").Value = -41260
$ws.Range("H61").Value = 7939093.5
$ws.Range("I61").Value = 8549739
$ws.Range("J61").Value = 698
$ws.Range("K61").Value = 8549739
$ws.Range("L61").Value = 698
$ws.Range("M61").Value = -8549527
$ws.Range("N61").Value = -1122
$ws.Range("H74").Value = 601.85297
$ws.Range("I74").Value = 501.68967
$ws.Range("J74").Value = 1182.8
$ws.Range("K74").Value = 501.68967
$ws.Range("L74").Value = 1182.8
$ws.Range("M74").Value = 372.31033
$ws.Range("N74").Value = -2930.8
$ws.Range("H77").Value = 601.85297
$ws.Range("I77").Value = 501.68967
$ws.Range("J77").Value = 1182.8
$ws.Range("K77").Value = 2508.44835
$ws.Range("L77").Value = 5914
$ws.Range("M77").Value = 1859.55165
$ws.Range("N77").Value = -14650
$ws.Range("H132").Value = 3061213.5
$ws.Range("I132").Value = 3685807.2
$ws.Range("J132").Value = 704.2
$ws.Range("K132").Value = 11057421.6
$ws.Range("L132").Value = 2112.6
$ws.Range("M132").Value = -11054891.6
$ws.Range("N132").Value = -7172.6
$ws.Range("H136").Value = 7939093.5
$ws.Range("I136").Value = 8549739
$ws.Range("J136").Value = 698
$ws.Range("K136").Value = 25649217
$ws.Range("L136").Value = 2094
$ws.Range("M136").Value = -25646667
$ws.Range("N136").Value = -7194
$ws.Range("H137").Value = 17200
$ws.Range("J137").Value = 17200
$ws.Range("L137").Value = 17200
$ws.Range("N137").Value = -27400

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21511382
$ws.Range("I31").Value = 25642086
$ws.Range("J31").Value = 31719.4
$ws.Range("K31").Value = 25642086
$ws.Range("L31").Value = 31719.4
$ws.Range("M31").Value = -25641791
$ws.Range("N31").Value = -32309.4
$ws.Range("H34").Value = 21511382
$ws.Range("I34").Value = 25642086
$ws.Range("J34").Value = 31719.4
$ws.Range("K34").Value = 25642086
$ws.Range("L34").Value = 31719.4
$ws.Range("M34").Value = -25641884
$ws.Range("N34").Value = -32123.4
$ws.Range("H62").Value = 3380.6
$ws.Range("J62").Value = 4301
$ws.Range("L62").Value = 4301
$ws.Range("N62").Value = -5549
$ws.Range("H65").Value = 3380.6
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 4301
$ws.Range("L65").Value = 21505
$ws.Range("N65").Value = -27745
$ws.Range("H98").Value = 38345
$ws.Range("J98").Value = 38345
$ws.Range("L98").Value = 38345
$ws.Range("N98").Value = -42837

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15933603
$ws.Range("I131").Value = 62500428
$ws.Range("J131").Value = 1605348.6
$ws.Range("K131").Value = 187501284
$ws.Range("L131").Value = 4816045.800000001
$ws.Range("M131").Value = -187496244
$ws.Range("N131").Value = -4826125.800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50081.816
$ws.Range("I70").Value = 71140
$ws.Range("J70").Value = 4957.143
$ws.Range("K70").Value = 71140
$ws.Range("L70").Value = 4957.143
$ws.Range("M70").Value = -70870
$ws.Range("N70").Value = -5497.143
$ws.Range("H73").Value = 50081.816
$ws.Range("I73").Value = 71140
$ws.Range("J73").Value = 4957.143
$ws.Range("K73").Value = 71140
$ws.Range("L73").Value = 4957.143
$ws.Range("M73").Value = -70204
$ws.Range("N73").Value = -6829.143
$ws.Range("H132").Value = 27419166
$ws.Range("I132").Value = 37699480
$ws.Range("J132").Value = 4990.7617
$ws.Range("K132").Value = 113098440
$ws.Range("L132").Value = 14972.2851
$ws.Range("M132").Value = -113095910
$ws.Range("N132").Value = -20032.2851

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2357.9666
$ws.Range("I82").Value = 1767.3158
$ws.Range("J82").Value = 3378.182
$ws.Range("K82").Value = 1767.3158
$ws.Range("L82").Value = 3378.182
$ws.Range("M82").Value = -1406.3158
$ws.Range("N82").Value = -4100.182
$ws.Range("H85").Value = 2357.9666
$ws.Range("I85").Value = 1767.3158
$ws.Range("J85").Value = 3378.182
$ws.Range("K85").Value = 1767.3158
$ws.Range("L85").Value = 3378.182
$ws.Range("M85").Value = -519.3158000000001
$ws.Range("N85").Value = -5874.182

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 19478570
$ws.Range("I132").Value = 25532842
$ws.Range("J132").Value = 12538308
$ws.Range("K132").Value = 76598526
$ws.Range("L132").Value = 37614924
$ws.Range("M132").Value = -76595996
$ws.Range("N132").Value = -37619984
